# update to manual status column;
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the manualStatus column (I) values for rows 11-13 from numeric 4 to text "[4]"
$ws.Range("I11").Value = "[4]"
$ws.Range("I12").Value = "[4]"
$ws.Range("I13").Value = "[4]"

# Adjust row heights for rows 12 and 13
$ws.Rows.Item(12).RowHeight = 13.8
$ws.Rows.Item(13).RowHeight = 13.8

# Widen column F (fastqFileName) to fit its content
# (the engine rounds ColumnWidth to the nearest 1/6-character pixel grid, so
# 52.17 is the closest input that lands exactly on the target stored width)
$ws.Columns.Item(6).ColumnWidth = 52.17

# Update the active cell selection
$ws.Range("I13").Select()
